$wb = $excel.ActiveWorkbook

$wsAdd   = $wb.Worksheets.Item(1)   # AddCustomerTest
$wsOpen  = $wb.Worksheets.Item(2)   # OpenAccountTest
$wsSuite = $wb.Worksheets.Item(3)   # testSuite

# --- testSuite: rename headers to lowercase and flip OpenAccountTest's runmode ---
$wsSuite.Range("A1").Value = "tcid"
$wsSuite.Range("B1").Value = "runmode"
$wsSuite.Range("B4").Value = "n"
$null = $wsSuite.Range("B4").Select()

# --- AddCustomerTest: flip the runmode for the Jorge/Souza test row ---
$wsAdd.Range("E5").Value = "n"

# Make AddCustomerTest the active/selected sheet & cell (matches activeTab=0)
[void]$wsAdd.Activate()
$null = $wsAdd.Range("F8").Select()
